{"js": "// Apply the RS-RI schema edits described in the commit diff.\n//\n// Strategy: locate each target table row by searching for a unique piece\n// of text that still exists in the \"before\" document, then navigate via\n// parentTableCell / parentRow / parentTable so the edit isn't dependent on\n// brittle fixed table/row indices.\nasync function findRow(searchText) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not locate text: \" + searchText);\n  }\n  const cell = results.items[0].parentTableCell;\n  const row = cell.parentRow;\n  row.cells.load(\"items\");\n  await context.sync();\n  return row;\n}\n\nfunction setCellText(row, cellIndex, text) {\n  const cell = row.cells.items[cellIndex];\n  cell.body.clear();\n  cell.body.insertText(text, Word.InsertLocation.start);\n}\n\n// --- \"mobilizedResources\" row: rename tag + bump cardinality ---\nconst mobResRow = await findRow(\"mobilizedResources\");\nsetCellText(mobResRow, 0, \"mobilizedResource\");\nsetCellText(mobResRow, 3, \"1..n\");\n\n// --- \"resourceID\" row: relabel + rewrite the description/format text ---\nconst resourceIdRow = await findRow(\"resourceID\");\nsetCellText(resourceIdRow, 1, \"ID Ressource partag\u00e9\");\nsetCellText(\n  resourceIdRow,\n  4,\n  \"ID partag\u00e9 unique de la ressource engag\u00e9e \\u000b\" +\n    \"{orgID}.R.{ID unique de la ressource partag\u00e9e}\\u000b\" +\n    \"Ou, uniquement dans le cas o\u00f9 un ID unique de ressource ne peut pas \u00eatre garanti par l'organisation propri\u00e9taire :\\u000b\" +\n    \"{orgID}.R.{ID du dossier partag\u00e9}.{num\u00e9ro d\u2019ordre chronologique ressource}\"\n);\n\n// --- \"missionID\" row: relabel the field name ---\n// Apply this edit *before* inserting the new row below, otherwise the row\n// reference captured by index could get reseated onto the freshly inserted\n// row instead of the original \"missionID\" row.\nconst missionIdRow = await findRow(\"missionID\");\nsetCellText(missionIdRow, 1, \"ID Mission local\");\n\n// Insert the new \"RSDRId\" row right after the resourceID row.\nresourceIdRow.insertRows(Word.InsertLocation.after, 1, [[\n  \"RSDRId\",\n  \"ID Demande partag\u00e9\",\n  \"string\",\n  \"0..1\",\n  \"Identifiant unique partag\u00e9 de la demande de ressource (si la ressource a \u00e9t\u00e9 engag\u00e9e suite \u00e0 une demande de ressource)\\u000b\" +\n    \"{orgID}.D.{ID unique de la demande dans le syst\u00e8me \u00e9metteur}\",\n  \"fr.health.samu770.D.1249875\"\n]]);\n\n// --- \"availability\" row: collapse the duplicated ENUM text into one line ---\nconst availabilityRow = await findRow(\"Disponibilit\u00e9 du vecteur\");\nsetCellText(availabilityRow, 2, \"string\\u000b(ENUM : DISPONIBLE, INDISPONIBLE, INCONNUE)\");\n\n// --- contact \"type\" row: extend the ENUM list ---\nconst contactTypeRow = await findRow(\"Type de contact\");\nsetCellText(contactTypeRow, 2, \"string\\u000b(ENUM: TEL, EMAIL, FAX, POSTAL, WEB, RADIO)\");\n\nawait context.sync();\n", "ps1": "# Apply the RS-RI schema edits described in the commit diff.\n#\n# Strategy: locate each target table row by searching for a unique piece of\n# text that still exists in the \"before\" document (via Find.Execute), then\n# resolve the owning table by matching range offsets against $d.Tables, and\n# finally index into that table's Rows collection. This avoids hard-coding\n# table/row numbers.\n\n$d = $word.ActiveDocument\n$nl = [char]11   # manual line break -> serializes as <w:br/>\n\nfunction Find-TableRow {\n  param($doc, [string]$searchText)\n\n  $range = $doc.Content.Duplicate\n  $found = $range.Find.Execute($searchText, $true)\n  if (-not $found) {\n    throw \"Could not locate text: $searchText\"\n  }\n\n  $targetTable = $null\n  for ($i = 1; $i -le $doc.Tables.Count; $i++) {\n    $t = $doc.Tables.Item($i)\n    if ($range.Start -ge $t.Range.Start -and $range.End -le $t.Range.End) {\n      $targetTable = $t\n      break\n    }\n  }\n  if ($null -eq $targetTable) {\n    throw \"Could not locate table containing: $searchText\"\n  }\n\n  $cell = $range.Cells.Item(1)\n  return $targetTable.Rows.Item($cell.RowIndex)\n}\n\n# --- \"mobilizedResources\" row: rename tag + bump cardinality ---\n$mobRow = Find-TableRow $d \"mobilizedResources\"\n$mobRow.Cells.Item(1).Range.Text = \"mobilizedResource\"\n$mobRow.Cells.Item(4).Range.Text = \"1..n\"\n\n# --- \"resourceID\" row: relabel + rewrite the description/format text ---\n$resourceRow = Find-TableRow $d \"resourceID\"\n$resourceRow.Cells.Item(2).Range.Text = \"ID Ressource partag\u00e9\"\n$resourceRow.Cells.Item(5).Range.Text = (\n  \"ID partag\u00e9 unique de la ressource engag\u00e9e \" + $nl +\n  \"{orgID}.R.{ID unique de la ressource partag\u00e9e}\" + $nl +\n  \"Ou, uniquement dans le cas o\u00f9 un ID unique de ressource ne peut pas \u00eatre garanti par l'organisation propri\u00e9taire :\" + $nl +\n  \"{orgID}.R.{ID du dossier partag\u00e9}.{num\u00e9ro d\u2019ordre chronologique ressource}\"\n)\n\n# --- \"missionID\" row: relabel the field name ---\n# Apply this edit *before* inserting the new row below, so the insertion\n# cannot disturb which physical row this refers to.\n$missionRow = Find-TableRow $d \"missionID\"\n$missionRow.Cells.Item(2).Range.Text = \"ID Mission local\"\n\n# Insert the new \"RSDRId\" row right before the missionID row (i.e. right\n# after the resourceID row).\n$resourceTable = $resourceRow.Range.Tables.Item(1)\n$newRow = $resourceTable.Rows.Add($missionRow)\n$newRow.Cells.Item(1).Range.Text = \"RSDRId\"\n$newRow.Cells.Item(2).Range.Text = \"ID Demande partag\u00e9\"\n$newRow.Cells.Item(3).Range.Text = \"string\"\n$newRow.Cells.Item(4).Range.Text = \"0..1\"\n$newRow.Cells.Item(5).Range.Text = (\n  \"Identifiant unique partag\u00e9 de la demande de ressource (si la ressource a \u00e9t\u00e9 engag\u00e9e suite \u00e0 une demande de ressource)\" + $nl +\n  \"{orgID}.D.{ID unique de la demande dans le syst\u00e8me \u00e9metteur}\"\n)\n$newRow.Cells.Item(6).Range.Text = \"fr.health.samu770.D.1249875\"\n\n# --- \"availability\" row: collapse the duplicated ENUM text into one line ---\n$availabilityRow = Find-TableRow $d \"Disponibilit\u00e9 du vecteur\"\n$availabilityRow.Cells.Item(3).Range.Text = \"string\" + $nl + \"(ENUM : DISPONIBLE, INDISPONIBLE, INCONNUE)\"\n\n# --- contact \"type\" row: extend the ENUM list ---\n$contactTypeRow = Find-TableRow $d \"Type de contact\"\n$contactTypeRow.Cells.Item(3).Range.Text = \"string\" + $nl + \"(ENUM: TEL, EMAIL, FAX, POSTAL, WEB, RADIO)\"\n"}
